# Fix js calendar. Fix Bozena and name order in resources
#
# - Normalize professor name order: "Iwaniec Joanna" -> "Joanna Iwaniec"
# - Fix misspelled professor name: "Giermek Bozena" -> "Giermek Bożena"
# - Minor column width adjustments left behind by the edit
# - Leave the cursor on the last touched cell (E76)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 76
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -eq "Iwaniec Joanna") {
        $cell.Value = "Joanna Iwaniec"
    } elseif ($val -eq "Giermek Bozena") {
        $cell.Value = "Giermek Bożena"
    }
}

# Column widths drift slightly (auto-adjustment side effect of the edits above)
$ws.Columns.Item(1).ColumnWidth = 20.3333333333333
$ws.Columns.Item(3).ColumnWidth = 11.3333333333333
$ws.Columns.Item(4).ColumnWidth = 36.8333333333333
$ws.Columns.Item(5).ColumnWidth = 21.8333333333333
$ws.Columns.Item(7).ColumnWidth = 12.5

# Selection ends on the last edited cell
[void]$ws.Range("E76").Select()
